$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.766.26'
$ws.Range("E2").Value = '  +2.57%  '
$ws.Range("D3").Value = '2.335.18'
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("E4").Value = '  -2.68%  '
$ws.Range("D5").Value = '''313.48'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '''108.42'
$ws.Range("E6").Value = '  +7.07%  '
$ws.Range("D7").Value = '''0.632'
$ws.Range("E7").Value = '  +1.87%  '
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").Value = '''0.621'
$ws.Range("E9").Value = '  +4.86%  '
$ws.Range("D10").Value = '''41.31'
$ws.Range("E10").Value = '  +7.58%  '
$ws.Range("E11").Value = '  +2.47%  '
$ws.Range("E12").Value = '  +4.68%  '
$ws.Range("E13").Value = '  +4.39%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '''15.49'
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("D16").Value = '2.686.35'
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").Value = '2.322.57'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '43.692.07'
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("E19").Value = '  +4.34%  '
$ws.Range("D21").Value = '''13.02'
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").Value = '''74.39'
$ws.Range("E22").Value = '  +2.29%  '
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").Value = '''268.45'
$ws.Range("E24").Value = '  +2.88%  '
$ws.Range("E25").Value = '  +5.71%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '''7.64'
$ws.Range("E27").Value = '  +12.31%  '
$ws.Range("E28").Value = '  +5.26%  '
$ws.Range("E29").Value = '  -1.01%  '
$ws.Range("D30").Value = '''39.39'
$ws.Range("E30").Value = '  +9.55%  '
$ws.Range("D31").Value = '''22.55'
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").Value = '''168.17'
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = '''0.0901'
$ws.Range("E33").Value = '  +4.83%  '
$ws.Range("E34").Value = '  +8.95%  '
$ws.Range("E35").Value = '  +1.75%  '
$ws.Range("E36").Value = '  +4.63%  '
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("E38").Value = '  +6.00%  '
$ws.Range("D39").Value = '''2.91'
$ws.Range("E39").Value = '  +10.89%  '
$ws.Range("E40").Value = '  +4.55%  '
$ws.Range("E41").Value = '  +11.14%  '
$ws.Range("D42").Value = '''103.74'
$ws.Range("E42").Value = '  +12.11%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.240'
$ws.Range("E43").Value = '  +6.53%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").Value = '''13.62'
$ws.Range("E44").Value = '  +15.09%  '
$ws.Range("D45").Value = '''71.80'
$ws.Range("E45").Value = '  +4.40%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '''114.85'
$ws.Range("E47").Value = '  +4.46%  '
$ws.Range("D48").Value = '''0.218'
$ws.Range("E48").Value = '  +18.92%  '
$ws.Range("D49").Value = '1.656.37'
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("D50").Value = '''9.01'
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("D51").Value = '''76.05'
$ws.Range("E51").Value = '  -2.88%  '
